# SummaryQuantifications.xlsx edit:
# Insert two new columns ("abun_con", "abun_inf") right after the "hours"
# column (i.e. before the old column C) and fill them with the new
# abundance-count data, then restore/clean up the formatting of the whole
# table (no cell borders anymore, consistent fonts/number formats) and move
# the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert two blank columns at C:D - this shifts the old C:K to E:M
# ---------------------------------------------------------------------
$ws.Columns("C:D").Insert()

# ---------------------------------------------------------------------
# 2. Headers for the two new columns
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 3).Value = "abun_con"
$ws.Cells.Item(1, 4).Value = "abun_inf"

# ---------------------------------------------------------------------
# 3. New column data (counts)
# ---------------------------------------------------------------------
$newData = @{
    2 = @(14782, 12877)
    3 = @(11302, 13857)
    4 = @(16589, 13971)
    5 = @(19439, 7719)
    6 = @(18097, 4155)
    7 = @(29564, 5325)
    8 = @(29266, 4686)
    9 = @(26167, 2563)
}
foreach ($r in $newData.Keys) {
    $vals = $newData[$r]
    $ws.Cells.Item($r, 3).Value = $vals[0]
    $ws.Cells.Item($r, 4).Value = $vals[1]
}

# ---------------------------------------------------------------------
# 4. Empty trailing column N (header + data rows + blank row 10)
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 14).Value = ""
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 14).Value = ""
}

# ---------------------------------------------------------------------
# 5. Extend the formatted-but-empty row 10 and row 11 to the new columns
# ---------------------------------------------------------------------
$ws.Range("K11:L11").Value = ""

Write-Host "edit applied"
